# Swap the pick data between "Julien Belvisi" (row 22) and "Michel Cusson" (row 26),
# and fill in Michel Cusson's previously-missing Round 2 (OUEST) picks.
#
# Commit message: "added missing picks round 2"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")

# --- Capture current ("before") values for the two rows we need to touch ---
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","T","U","V","W","X","Y","Z","AA")

$row22 = @{}
$row26 = @{}
foreach ($c in $cols) {
    $row22[$c] = $ws.Range("$c" + "22").Value()
    $row26[$c] = $ws.Range("$c" + "26").Value()
}

# --- Swap the person names (shared strings 70/71) between the rows ---
$name22 = $ws.Range("B22").Value()
$name26 = $ws.Range("B26").Value()
$ws.Range("B22").Value = $name26
$ws.Range("B26").Value = $name22

# --- Swap Round 1 + Round 2 (EST + OUEST) picks between the two rows ---
foreach ($c in $cols) {
    $ws.Range("$c" + "22").Value = $row26[$c]
    $ws.Range("$c" + "26").Value = $row22[$c]
}

# --- Fill in Michel Cusson's missing Round 2 OUEST picks (now on row 22) ---
$ws.Range("Y22").Value = 5
$ws.Range("Z22").Value = 7

# --- Update the active selection as recorded in the workbook view ---
$ws.Range("B5").Select()
